$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row columns: "<name>_old" -> "<name>_FV2410" and
#        "<name>_new" -> "<name>_FV2504" (column "diff" stays untouched). ---
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $colOld = $i + 1
    $ws.Cells.Item(1, $colOld).Value = "$($baseNames[$i])_FV2410"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $colNew = $i + 12
    $ws.Cells.Item(1, $colNew).Value = "$($baseNames[$i])_FV2504"
}

# --- 2. Freeze the header row (split after row 1). ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table ("Table1"). ---
$dataRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

Write-Host "Headers renamed, panes frozen, table '$($tbl.Name)' created."
